# Auto-generated update: refresh "剩余" (remaining days) and "开始时间" (start date)
# columns for each data row based on a new reference date (2026-01-07).
# For rows where the remaining-day countdown would reach zero or below,
# the cycle restarts: start date (F) is set to the new reference date and
# remaining (E) is reset to the full duration (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; E=5; F=20251229},
    @{Row=3; E=5; F=20251229},
    @{Row=4; E=5; F=20251229},
    @{Row=5; E=7; F=20260104},
    @{Row=6; E=5; F=20251229},
    @{Row=7; E=7; F=20260104},
    @{Row=8; E=5; F=20251229},
    @{Row=9; E=7; F=20260104},
    @{Row=10; E=5; F=20260105},
    @{Row=11; E=5; F=20251229},
    @{Row=12; E=7; F=20260104},
    @{Row=13; E=5; F=20251229},
    @{Row=14; E=5; F=20251229},
    @{Row=15; E=5; F=20251229},
    @{Row=16; E=1; F=20251229},
    @{Row=17; E=7; F=20260104},
    @{Row=18; E=10; F=20260107},
    @{Row=19; E=10; F=20260107},
    @{Row=20; E=10; F=20260107},
    @{Row=21; E=10; F=20260107},
    @{Row=22; E=7; F=20260104},
    @{Row=23; E=7; F=20260104},
    @{Row=24; E=7; F=20260104},
    @{Row=25; E=7; F=20260104},
    @{Row=26; E=7; F=20260104},
    @{Row=27; E=6; F=20260106},
    @{Row=28; E=10; F=20260107},
    @{Row=29; E=10; F=20260107},
    @{Row=30; E=10; F=20260107},
    @{Row=31; E=10; F=20260107},
    @{Row=32; E=10; F=20260107},
    @{Row=33; E=10; F=20260107},
    @{Row=34; E=10; F=20260107},
    @{Row=35; E=10; F=20260107},
    @{Row=37; E=10; F=20260107},
    @{Row=38; E=10; F=20260107},
    @{Row=39; E=10; F=20260107},
    @{Row=40; E=5; F=20260105},
    @{Row=41; E=5; F=20260105},
    @{Row=42; E=10; F=20260107},
    @{Row=43; E=7; F=20260104},
    @{Row=44; E=5; F=20260105},
    @{Row=45; E=7; F=20260104},
    @{Row=46; E=5; F=20260105},
    @{Row=47; E=10; F=20260107},
    @{Row=48; E=5; F=20260105},
    @{Row=49; E=6; F=20260106},
    @{Row=50; E=5; F=20260102},
    @{Row=51; E=5; F=20260102},
    @{Row=52; E=5; F=20260102},
    @{Row=53; E=5; F=20260102},
    @{Row=54; E=5; F=20260102},
    @{Row=55; E=5; F=20260102},
    @{Row=56; E=5; F=20260102},
    @{Row=57; E=5; F=20260102},
    @{Row=58; E=9; F=20260106},
    @{Row=59; E=9; F=20260106},
    @{Row=60; E=9; F=20260106},
    @{Row=61; E=6; F=20260106},
    @{Row=62; E=9; F=20260106},
    @{Row=63; E=9; F=20260106},
    @{Row=64; E=9; F=20260106},
    @{Row=65; E=10; F=20260107},
    @{Row=66; E=10; F=20260107},
    @{Row=67; E=10; F=20260107},
    @{Row=68; E=10; F=20260107},
    @{Row=69; E=10; F=20260107},
    @{Row=70; E=1; F=20251229},
    @{Row=71; E=1; F=20251229},
    @{Row=72; E=1; F=20251229},
    @{Row=73; E=1; F=20251229},
    @{Row=74; E=1; F=20251229},
    @{Row=75; E=1; F=20251229},
    @{Row=76; E=1; F=20251229},
    @{Row=77; E=4; F=20260101},
    @{Row=78; E=4; F=20260101},
    @{Row=79; E=4; F=20260101},
    @{Row=80; E=4; F=20260101},
    @{Row=81; E=4; F=20260101},
    @{Row=82; E=4; F=20260101},
    @{Row=83; E=4; F=20260101},
    @{Row=84; E=4; F=20260101},
    @{Row=85; E=4; F=20260101},
    @{Row=86; E=4; F=20260101},
    @{Row=87; E=5; F=20260105},
    @{Row=88; E=5; F=20260105},
    @{Row=89; E=5; F=20260105},
    @{Row=90; E=5; F=20260105},
    @{Row=91; E=7; F=20260104},
    @{Row=92; E=5; F=20260105},
    @{Row=93; E=4; F=20260101},
    @{Row=94; E=1; F=20260101},
    @{Row=95; E=3; F=20251231},
    @{Row=96; E=1; F=20251229},
    @{Row=97; E=1; F=20251229},
    @{Row=98; E=1; F=20251229},
    @{Row=99; E=1; F=20251229}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}
